$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'34.011.30"
$ws.Range("E2").Value = "  +3.20%  "

$ws.Range("D3").Value = "'1.788.65"
$ws.Range("E3").Value = "  +1.82%  "

$ws.Range("E4").Value = "  +0.51%  "

$ws.Range("D5").Value = "'226.38"
$ws.Range("E5").Value = "  -0.22%  "

$ws.Range("D6").Value = "'0.560"
$ws.Range("E6").Value = "  +3.04%  "

$ws.Range("E7").Value = "  +0.70%  "

$ws.Range("D8").Value = "'30.28"
$ws.Range("E8").Value = "  -5.17%  "

$ws.Range("D9").Value = "'46.72"
$ws.Range("E9").Value = "  +3.84%  "

$ws.Range("E10").Value = "  +0.41%  "

$ws.Range("E11").Value = "  +0.67%  "

$ws.Range("D12").Value = "'0.0925"
$ws.Range("E12").Value = "  +0.81%  "

$ws.Range("E13").Value = "  +1.88%  "

$ws.Range("D14").Value = "'1.783.03"
$ws.Range("E14").Value = "  +1.15%  "

$ws.Range("E15").Value = "  -0.83%  "

$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "'10.45"
$ws.Range("E16").Value = "  -0.59%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "'34.023.21"
$ws.Range("E17").Value = "  +3.31%  "

$ws.Range("E18").Value = "  -2.05%  "

$ws.Range("D19").Value = "'69.20"
$ws.Range("E19").Value = "  +0.65%  "

$ws.Range("D20").Value = "'252.34"
$ws.Range("E20").Value = "  -2.76%  "

$ws.Range("D21").Value = "'0.0₃0742"
$ws.Range("E21").Value = "  +0.27%  "

$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  +0.19%  "

$ws.Range("D23").Value = "'10.36"
$ws.Range("E23").Value = "  -1.28%  "

$ws.Range("D24").Value = "'4.23"
$ws.Range("E24").Value = "  -2.75%  "

$ws.Range("D25").Value = "'2.13"
$ws.Range("E25").Value = "  -2.00%  "

$ws.Range("D26").Value = "'158.31"
$ws.Range("E26").Value = "  -1.08%  "

$ws.Range("D27").Value = "'16.52"
$ws.Range("E27").Value = "  -0.12%  "

$ws.Range("E28").Value = "  -0.15%  "

$ws.Range("D29").Value = "'7.01"
$ws.Range("E29").Value = "  +0.53%  "

$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.61%  "

$ws.Range("E31").Value = "  -1.64%  "

$ws.Range("E32").Value = "  -0.08%  "

$ws.Range("E33").Value = "  +1.48%  "

$ws.Range("D34").Value = "'3.59"
$ws.Range("E34").Value = "  +3.42%  "

$ws.Range("E35").Value = "  +3.90%  "

$ws.Range("D36").Value = "'1.506.38"
$ws.Range("E36").Value = "  -3.01%  "

$ws.Range("E37").Value = "  +2.12%  "

$ws.Range("D38").Value = "'0.636"
$ws.Range("E38").Value = "  +0.60%  "

$ws.Range("D39").Value = "'0.0186"
$ws.Range("E39").Value = "  -0.01%  "

$ws.Range("D40").Value = "'83.62"
$ws.Range("E40").Value = "  -1.56%  "

$ws.Range("B41").Value = "HuobiToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D41").Value = "'2.36"
$ws.Range("E41").Value = "  +2.57%  "

$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").Value = "'2.71"
$ws.Range("E42").Value = "  -0.53%  "

$ws.Range("D43").Value = "'0.903"
$ws.Range("E43").Value = "  +3.25%  "

$ws.Range("D44").Value = "'0.0518"
$ws.Range("E44").Value = "  +0.69%  "

$ws.Range("E45").Value = "  -2.17%  "

$ws.Range("E46").Value = "  +2.16%  "

$ws.Range("D47").Value = "'1.947.15"
$ws.Range("E47").Value = "  +2.18%  "

$ws.Range("D48").Value = "'5.75"
$ws.Range("E48").Value = "  +1.02%  "

$ws.Range("E49").Value = "  +0.61%  "

$ws.Range("D50").Value = "'11.80"
$ws.Range("E50").Value = "  +6.09%  "

$ws.Range("D51").Value = "'51.53"
$ws.Range("E51").Value = "  -6.37%  "
